# Update countries & provincias Spain
# Applies the 24-Abril-2020 data refresh (11:52 -> 12:22) to the "Ciudades" sheet:
#  - a handful of place names shift row (new leaders overtake old ones in the
#    sort-by-cases list), so the label in column A for a few rows changes
#  - the Casos totales / Casos activos / Recuperados / Muertes columns (B:E)
#    get refreshed numbers for the affected rows
#  - the "Datos actualizados..." timestamp banner in A1 is bumped

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# --- timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 12:22"

# --- province/city names that moved row as the ranking reshuffled ---
$ws.Range("A14").Value = "Navarra"
$ws.Range("A15").Value = "Valencia/Valencia"
$ws.Range("A16").Value = "La Rioja"
$ws.Range("A17").Value = "Toledo"
$ws.Range("A31").Value = "Cantabria"
$ws.Range("A32").Value = "Caceres"

# --- refreshed case numbers (Casos totales, Casos activos, Recuperados, Muertes) ---
$ws.Range("B4").Value = 61726
$ws.Range("C4").Value = 34212
$ws.Range("D4").Value = 19749
$ws.Range("E4").Value = 7765

$ws.Range("B5").Value = 46571
$ws.Range("C5").Value = 16138
$ws.Range("D5").Value = 26040
$ws.Range("E5").Value = 4393

$ws.Range("B6").Value = 18053
$ws.Range("C6").Value = 4782
$ws.Range("D6").Value = 11016
$ws.Range("E6").Value = 2255

$ws.Range("B7").Value = 17776
$ws.Range("C7").Value = 5943
$ws.Range("D7").Value = 10221
$ws.Range("E7").Value = 1612

$ws.Range("B8").Value = 13780
$ws.Range("C8").Value = 8459
$ws.Range("D8").Value = 4128
$ws.Range("E8").Value = 1193

$ws.Range("B9").Value = 12495
$ws.Range("C9").Value = 3992
$ws.Range("D9").Value = 7396
$ws.Range("E9").Value = 1107

$ws.Range("B10").Value = 8932
$ws.Range("C10").Value = 1730
$ws.Range("D10").Value = 6820
$ws.Range("E10").Value = 382

$ws.Range("B13").Value = 5295
$ws.Range("C13").Value = 1781
$ws.Range("D13").Value = 2819
$ws.Range("E13").Value = 695

$ws.Range("B14").Value = 5180
$ws.Range("C14").Value = 1552
$ws.Range("D14").Value = 3211
$ws.Range("E14").Value = 417

$ws.Range("B15").Value = 5131
$ws.Range("C15").Value = 2194
$ws.Range("D15").Value = 2583
$ws.Range("E15").Value = 515

$ws.Range("B16").Value = 4865
$ws.Range("C16").Value = 1964
$ws.Range("D16").Value = 2594
$ws.Range("E16").Value = 307

$ws.Range("B17").Value = 3938
$ws.Range("C17").Value = 4178
$ws.Range("D17").Value = 10597
$ws.Range("E17").Value = 504

$ws.Range("B21").Value = 3403
$ws.Range("C21").Value = 1358
$ws.Range("D21").Value = 1634
$ws.Range("E21").Value = 411

$ws.Range("B26").Value = 2509
$ws.Range("C26").Value = 688
$ws.Range("D26").Value = 1590
$ws.Range("E26").Value = 231

$ws.Range("B31").Value = 2273
$ws.Range("C31").Value = 913
$ws.Range("D31").Value = 1182
$ws.Range("E31").Value = 178

$ws.Range("B32").Value = 2220
$ws.Range("C32").Value = 422
$ws.Range("D32").Value = 1482
$ws.Range("E32").Value = 316

$ws.Range("B33").Value = 2140
$ws.Range("C33").Value = 1017
$ws.Range("D33").Value = 995
$ws.Range("E33").Value = 128

$ws.Range("B36").Value = 1741
$ws.Range("C36").Value = 791
$ws.Range("D36").Value = 825
$ws.Range("E36").Value = 125

$ws.Range("B58").Value = 125
$ws.Range("C58").Value = 84
$ws.Range("D58").Value = 37

$ws.Range("B59").Value = 116
$ws.Range("C59").Value = 61
$ws.Range("D59").Value = 53
